$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.228.14"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.854.65"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'0.7008"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "'237.62"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "'0.07909"
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("D9").Value = "'0.3015"
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").Value = "'23.56"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").Value = "'0.08181"
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("D12").Value = "1.847.22"
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("D13").Value = "'5.181"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").Value = "'0.7039"
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").Value = "'89.47"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "29.221.69"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "'5.805"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").Value = "'0.000007814"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").Value = "'236.16"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "2.092.32"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").Value = "'7.495"
$ws.Range("D25").Value = "'162.82"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").Value = "'8.845"
$ws.Range("E26").Value = "  -1.58%  "
$ws.Range("D27").Value = "'0.1415"
$ws.Range("E27").Value = "  -2.21%  "
$ws.Range("D28").Value = "'18.05"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "'1.910"
$ws.Range("E29").Value = "  -2.74%  "
$ws.Range("D30").Value = "'1.407"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "'1.470"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").Value = "'4.316"
$ws.Range("E32").Value = "  -4.62%  "
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").Value = "'0.05151"
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("D35").Value = "'1.162"
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("D36").Value = "'0.7095"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").Value = "'0.9950"
$ws.Range("E37").Value = "  -2.72%  "
$ws.Range("D38").Value = "'2.680"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").Value = "'0.01844"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").Value = "'2.708"
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("D41").Value = "1.155.66"
$ws.Range("E41").Value = "  +4.82%  "
$ws.Range("D42").Value = "'0.9287"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").Value = "'5.975"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").Value = "'0.4239"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("D45").Value = "'70.09"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "'102.59"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "'0.5290"
$ws.Range("E48").Value = "  -2.90%  "
$ws.Range("D49").Value = "'1.735"
$ws.Range("E49").Value = "  -3.21%  "
$ws.Range("D50").Value = "'9.119"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").Value = "'6.945"
$ws.Range("E51").Value = "  -1.13%  "

Write-Host "Applied cryptos update."
